$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append two new rows (46-47) ---
$proximity = $wb.Worksheets.Item("Proximity")

# Ensure the Date column is stored as text (matching the existing log rows)
# rather than being auto-converted into a date serial number.
$proximity.Range("A46:A47").NumberFormat = "@"

$proximity.Range("A46").Value = "2026-02-01"
$proximity.Range("B46").Value = "14:45:42"
$proximity.Range("C46").Value = "14:00"
$proximity.Range("D46").Value = "Living Room Main Door"
$proximity.Range("E46").Value = "ENTER"
$proximity.Range("F46").Value = "User ENTERED Living Room Main Door"

$proximity.Range("A47").Value = "2026-02-01"
$proximity.Range("B47").Value = "14:45:45"
$proximity.Range("C47").Value = "14:00"
$proximity.Range("D47").Value = "Living Room Main Door"
$proximity.Range("E47").Value = "EXIT"
$proximity.Range("F47").Value = "User EXITED Living Room Main Door"

# Strip the leftover explicit style index so the new cells fall back to the
# default (unstyled) formatting used throughout the rest of the log.
$proximity.Range("A46:F47").ClearFormats()

# --- Camera sheet: append two new rows (31-32) ---
$camera = $wb.Worksheets.Item("Camera")

$camera.Range("A31:A32").NumberFormat = "@"

$camera.Range("A31").Value = "2026-02-01"
$camera.Range("B31").Value = "14:45:44"
$camera.Range("C31").Value = "14:00"
$camera.Range("D31").Value = "Living Room Main Door"
$camera.Range("E31").Value = "Image Captured"
$camera.Range("F31").Value = "Active"

$camera.Range("A32").Value = "2026-02-01"
$camera.Range("B32").Value = "14:45:45"
$camera.Range("C32").Value = "14:00"
$camera.Range("D32").Value = "Living Room Main Door"
$camera.Range("E32").Value = "Image Received"
$camera.Range("F32").Value = "Active"

$camera.Range("A31:F32").ClearFormats()
